$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly update swaps the "Primera" quality pair (rows 2 & 4) and the
# "Segunda" quality pair (rows 3 & 5): row 2 now carries what used to be in
# row 4, and row 4 now carries what used to be in row 2 (likewise for 3/5).
# Only the date (D) and price columns (N, O, P, S) change; everything else
# stays put.

# Capture the original values for the two rows in each swapped pair first,
# so overwriting one row doesn't clobber the source data for the other.
$row2 = @{
    D = $ws.Range("D2").Value2
    N = $ws.Range("N2").Value2
    O = $ws.Range("O2").Value2
    P = $ws.Range("P2").Value2
    S = $ws.Range("S2").Value2
}
$row4 = @{
    D = $ws.Range("D4").Value2
    N = $ws.Range("N4").Value2
    O = $ws.Range("O4").Value2
    P = $ws.Range("P4").Value2
    S = $ws.Range("S4").Value2
}
$row3 = @{
    D = $ws.Range("D3").Value2
    N = $ws.Range("N3").Value2
    O = $ws.Range("O3").Value2
    P = $ws.Range("P3").Value2
    S = $ws.Range("S3").Value2
}
$row5 = @{
    D = $ws.Range("D5").Value2
    N = $ws.Range("N5").Value2
    O = $ws.Range("O5").Value2
    P = $ws.Range("P5").Value2
    S = $ws.Range("S5").Value2
}

# Row 2 <- old row 4 values
$ws.Range("D2").Value = $row4.D
$ws.Range("N2").Value = $row4.N
$ws.Range("O2").Value = $row4.O
$ws.Range("P2").Value = $row4.P
$ws.Range("S2").Value = $row4.S

# Row 4 <- old row 2 values
$ws.Range("D4").Value = $row2.D
$ws.Range("N4").Value = $row2.N
$ws.Range("O4").Value = $row2.O
$ws.Range("P4").Value = $row2.P
$ws.Range("S4").Value = $row2.S

# Row 3 <- old row 5 values
$ws.Range("D3").Value = $row5.D
$ws.Range("N3").Value = $row5.N
$ws.Range("O3").Value = $row5.O
$ws.Range("P3").Value = $row5.P
$ws.Range("S3").Value = $row5.S

# Row 5 <- old row 3 values
$ws.Range("D5").Value = $row3.D
$ws.Range("N5").Value = $row3.N
$ws.Range("O5").Value = $row3.O
$ws.Range("P5").Value = $row3.P
$ws.Range("S5").Value = $row3.S
